$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M4").Value = -108.0
$ws.Range("H4").Value = 222.0
$ws.Range("I4").Value = 222.0
$ws.Range("K4").Value = 222.0
$ws.Range("J17").Value = 2568.9092
$ws.Range("L17").Value = 7706.7276
$ws.Range("N17").Value = -8042.7276
$ws.Range("H17").Value = 2568.9092
$ws.Range("I28").Value = 177.25
$ws.Range("K28").Value = 177.25
$ws.Range("H28").Value = 341.8
$ws.Range("M28").Value = 307.75
$ws.Range("K31").Value = 1198.5
$ws.Range("H31").Value = 399.5
$ws.Range("I31").Value = 399.5
$ws.Range("M31").Value = -968.5
$ws.Range("J41").Value = 877.0
$ws.Range("H41").Value = 786.0
$ws.Range("L41").Value = 877.0
$ws.Range("N41").Value = -1757.0
$ws.Range("L45").Value = 30000.0
$ws.Range("N45").Value = -30384.0
$ws.Range("H45").Value = 3760.0
$ws.Range("J45").Value = 10000.0
$ws.Range("I74").Value = 8623.75
$ws.Range("K74").Value = 8623.75
$ws.Range("H74").Value = 10115.571
$ws.Range("M74").Value = -7687.75
$ws.Range("H77").Value = 10115.571
$ws.Range("K77").Value = 43118.75
$ws.Range("I77").Value = 8623.75
$ws.Range("M77").Value = -38438.75
$ws.Range("K80").Value = 2160.9999
$ws.Range("I80").Value = 720.3333
$ws.Range("M80").Value = -1162.9999
$ws.Range("H80").Value = 848.3
$ws.Range("M83").Value = -1490.9997
$ws.Range("H83").Value = 848.3
$ws.Range("K83").Value = 6482.9997
$ws.Range("I83").Value = 720.3333
$ws.Range("H112").Value = 4339.0835
$ws.Range("N112").Value = -15484.4858
$ws.Range("L112").Value = 13268.4858
$ws.Range("J112").Value = 4422.8286
$ws.Range("H132").Value = 2509.775
$ws.Range("K132").Value = 6337.769400000001
$ws.Range("M132").Value = -3807.769400000001
$ws.Range("I132").Value = 2112.5898
$ws.Range("H137").Value = 9034.286
$ws.Range("I137").Value = 10106.667
$ws.Range("M137").Value = -27770.001
$ws.Range("K137").Value = 30320.001
$ws.Range("K138").Value = 6093.1362
$ws.Range("I138").Value = 2031.0454
$ws.Range("H138").Value = 3330.5518
$ws.Range("M138").Value = -953.1361999999999
$ws.Range("J141").Value = 66499.5
$ws.Range("H141").Value = 8535.913
$ws.Range("I141").Value = 3015.5715
$ws.Range("L141").Value = 199498.5
$ws.Range("K141").Value = 9046.7145
$ws.Range("M141").Value = -3866.7145
$ws.Range("N141").Value = -209858.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L2").Value = 1027.5
$ws.Range("H2").Value = 765.63635
$ws.Range("J2").Value = 1027.5
$ws.Range("N2").Value = -1253.5
$ws.Range("K2").Value = 667.4375
$ws.Range("I2").Value = 667.4375
$ws.Range("M2").Value = -554.4375
$ws.Range("M32").Value = -3091.5
$ws.Range("I32").Value = 3378.5
$ws.Range("H32").Value = 4752.0
$ws.Range("K32").Value = 3378.5
$ws.Range("M45").Value = -1916.5
$ws.Range("I45").Value = 2293.5
$ws.Range("K45").Value = 2293.5
$ws.Range("H45").Value = 2109.2
$ws.Range("H61").Value = 2803.1428
$ws.Range("M61").Value = -2160.606
$ws.Range("I61").Value = 2372.606
$ws.Range("K61").Value = 2372.606
$ws.Range("H63").Value = 3145.3333
$ws.Range("I63").Value = 3509.25
$ws.Range("K63").Value = 3509.25
$ws.Range("N63").Value = -3789.5
$ws.Range("M63").Value = -2823.25
$ws.Range("J63").Value = 2417.5
$ws.Range("L63").Value = 2417.5
$ws.Range("L66").Value = 12087.5
$ws.Range("J66").Value = 2417.5
$ws.Range("I66").Value = 3509.25
$ws.Range("H66").Value = 3145.3333
$ws.Range("K66").Value = 17546.25
$ws.Range("N66").Value = -18951.5
$ws.Range("M66").Value = -14114.25
$ws.Range("I74").Value = 5825.0586
$ws.Range("K74").Value = 5825.0586
$ws.Range("H74").Value = 9124.7
$ws.Range("M74").Value = -4951.0586
$ws.Range("H77").Value = 9124.7
$ws.Range("K77").Value = 29125.293
$ws.Range("I77").Value = 5825.0586
$ws.Range("M77").Value = -24757.293
$ws.Range("I88").Value = 905.0909
$ws.Range("L88").Value = 788.3333
$ws.Range("M88").Value = -499.0909
$ws.Range("N88").Value = -1600.3333
$ws.Range("H88").Value = 852.55
$ws.Range("J88").Value = 788.3333
$ws.Range("K88").Value = 905.0909
$ws.Range("M91").Value = 498.9091
$ws.Range("K91").Value = 905.0909
$ws.Range("N91").Value = -3596.3333
$ws.Range("L91").Value = 788.3333
$ws.Range("H91").Value = 852.55
$ws.Range("J91").Value = 788.3333
$ws.Range("I91").Value = 905.0909
$ws.Range("K97").Value = 547.36365
$ws.Range("M97").Value = -51.36365000000001
$ws.Range("I97").Value = 547.36365
$ws.Range("H97").Value = 532.4167
$ws.Range("H116").Value = 765.63635
$ws.Range("M116").Value = 1626.5625
$ws.Range("J116").Value = 1027.5
$ws.Range("L116").Value = 1027.5
$ws.Range("N116").Value = -5615.5
$ws.Range("I116").Value = 667.4375
$ws.Range("K116").Value = 667.4375
$ws.Range("K122").Value = 13234.5
$ws.Range("I122").Value = 4411.5
$ws.Range("M122").Value = -10784.5
$ws.Range("H122").Value = 4396.2593
$ws.Range("H136").Value = 2803.1428
$ws.Range("I136").Value = 2372.606
$ws.Range("K136").Value = 7117.818000000001
$ws.Range("M136").Value = -4567.818000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("N3").Value = -1255.5
$ws.Range("J3").Value = 1027.5
$ws.Range("M3").Value = -553.4375
$ws.Range("K3").Value = 667.4375
$ws.Range("I3").Value = 667.4375
$ws.Range("L3").Value = 1027.5
$ws.Range("H3").Value = 765.63635
$ws.Range("I20").Value = 1888.3077
$ws.Range("K20").Value = 1888.3077
$ws.Range("M20").Value = -1641.3077
$ws.Range("H20").Value = 2442.7896
$ws.Range("H94").Value = 1758.0
$ws.Range("I94").Value = 1083.8718
$ws.Range("K94").Value = 1083.8718
$ws.Range("M94").Value = -632.8717999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 108.947365
$ws.Range("M7").Value = 21.0
$ws.Range("N7").Value = -364.0
$ws.Range("I7").Value = 92.0
$ws.Range("J7").Value = 138.0
$ws.Range("K7").Value = 92.0
$ws.Range("L7").Value = 138.0
$ws.Range("J22").Value = 1267.3077
$ws.Range("L22").Value = 1267.3077
$ws.Range("N22").Value = -1967.3077
$ws.Range("H22").Value = 834.0
$ws.Range("L31").Value = 880816.7
$ws.Range("H31").Value = 547044.1
$ws.Range("N31").Value = -881406.7
$ws.Range("J31").Value = 880816.7
$ws.Range("H34").Value = 547044.1
$ws.Range("L34").Value = 880816.7
$ws.Range("J34").Value = 880816.7
$ws.Range("N34").Value = -881220.7
$ws.Range("L62").Value = 4377.2
$ws.Range("N62").Value = -5625.2
$ws.Range("K62").Value = 3322.7222
$ws.Range("J62").Value = 4377.2
$ws.Range("I62").Value = 3322.7222
$ws.Range("H62").Value = 3551.9565
$ws.Range("M62").Value = -2698.7222
$ws.Range("J65").Value = 4377.2
$ws.Range("N65").Value = -28126.0
$ws.Range("M65").Value = -13493.611
$ws.Range("I65").Value = 3322.7222
$ws.Range("L65").Value = 21886.0
$ws.Range("K65").Value = 16613.611
$ws.Range("H65").Value = 3551.9565
$ws.Range("L133").Value = 44400.0
$ws.Range("H133").Value = 44400.0
$ws.Range("J133").Value = 44400.0
$ws.Range("N133").Value = -49460.0
$ws.Range("I134").Value = 4420.3335
$ws.Range("K134").Value = 13261.0005
$ws.Range("H134").Value = 3072.2
$ws.Range("M134").Value = -10726.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M3").Value = -1369.57145
$ws.Range("K3").Value = 1481.57145
$ws.Range("I3").Value = 493.85715
$ws.Range("H3").Value = 682.125
$ws.Range("I11").Value = 459.4
$ws.Range("H11").Value = 632.8333
$ws.Range("K11").Value = 1378.2
$ws.Range("M11").Value = -1238.2
$ws.Range("K56").Value = 16566.75
$ws.Range("H56").Value = 16566.75
$ws.Range("I56").Value = 16566.75
$ws.Range("M56").Value = -16036.75
$ws.Range("H86").Value = 0.0
$ws.Range("L86").Value = 0.0
$ws.Range("J86").Value = 0.0
$ws.Range("N86").Value = $null
$ws.Range("H89").Value = 0.0
$ws.Range("L89").Value = 0.0
$ws.Range("J89").Value = 0.0
$ws.Range("N89").Value = $null
$ws.Range("N106").Value = -36890.999
$ws.Range("L106").Value = 34998.999
$ws.Range("H106").Value = 12141.143
$ws.Range("J106").Value = 11666.333
$ws.Range("L107").Value = 6000.0
$ws.Range("J107").Value = 2000.0
$ws.Range("N107").Value = -9840.0
$ws.Range("H107").Value = 1592.5834
$ws.Range("N131").Value = -16162.125
$ws.Range("I131").Value = 799.3333
$ws.Range("H131").Value = 1781.7667
$ws.Range("M131").Value = 2642.0001
$ws.Range("K131").Value = 2397.9999
$ws.Range("L131").Value = 6082.125
$ws.Range("J131").Value = 2027.375
$ws.Range("I140").Value = 1117.3478
$ws.Range("M140").Value = 1827.9566
$ws.Range("K140").Value = 3352.0434
$ws.Range("H140").Value = 1117.3478

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 23422.516
$ws.Range("M70").Value = -34639.117
$ws.Range("J70").Value = 11218.0
$ws.Range("K70").Value = 34909.117
$ws.Range("N70").Value = -11758.0
$ws.Range("L70").Value = 11218.0
$ws.Range("I70").Value = 34909.117
$ws.Range("I73").Value = 34909.117
$ws.Range("J73").Value = 11218.0
$ws.Range("K73").Value = 34909.117
$ws.Range("L73").Value = 11218.0
$ws.Range("N73").Value = -13090.0
$ws.Range("H73").Value = 23422.516
$ws.Range("M73").Value = -33973.117
$ws.Range("K80").Value = 5884.8335
$ws.Range("I80").Value = 5884.8335
$ws.Range("J80").Value = 2325.0
$ws.Range("M80").Value = -4886.8335
$ws.Range("N80").Value = -4321.0
$ws.Range("H80").Value = 4104.9165
$ws.Range("L80").Value = 2325.0
$ws.Range("M83").Value = -24432.1675
$ws.Range("H83").Value = 4104.9165
$ws.Range("J83").Value = 2325.0
$ws.Range("K83").Value = 29424.1675
$ws.Range("N83").Value = -21609.0
$ws.Range("I83").Value = 5884.8335
$ws.Range("L83").Value = 11625.0
$ws.Range("H94").Value = 47498.5
$ws.Range("L94").Value = 47498.5
$ws.Range("N94").Value = -48850.5
$ws.Range("J94").Value = 47498.5
$ws.Range("M99").Value = -2751.5
$ws.Range("K99").Value = 4997.5
$ws.Range("H99").Value = 19998.0
$ws.Range("I99").Value = 4997.5
$ws.Range("L124").Value = 70000.0
$ws.Range("N124").Value = -79820.0
$ws.Range("H124").Value = 70000.0
$ws.Range("J124").Value = 70000.0
$ws.Range("I126").Value = 5722.4
$ws.Range("K126").Value = 17167.2
$ws.Range("N126").Value = -34924.001
$ws.Range("J126").Value = 9994.667
$ws.Range("H126").Value = 7324.5
$ws.Range("M126").Value = -14697.2
$ws.Range("L126").Value = 29984.001
$ws.Range("H132").Value = 4628.3555
$ws.Range("K132").Value = 10463.5716
$ws.Range("M132").Value = -7933.571599999999
$ws.Range("I132").Value = 3487.8572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J22").Value = 1055.6809
$ws.Range("L22").Value = 1055.6809
$ws.Range("I22").Value = 125846.625
$ws.Range("N22").Value = -1645.6809
$ws.Range("M22").Value = -125551.625
$ws.Range("K22").Value = 125846.625
$ws.Range("H22").Value = 19207.092
$ws.Range("H27").Value = 19207.092
$ws.Range("I27").Value = 125846.625
$ws.Range("L27").Value = 1055.6809
$ws.Range("M27").Value = -125739.625
$ws.Range("J27").Value = 1055.6809
$ws.Range("K27").Value = 125846.625
$ws.Range("N27").Value = -1269.6809
$ws.Range("J68").Value = 6074.875
$ws.Range("M68").Value = -1075.6666
$ws.Range("I68").Value = 1824.6666
$ws.Range("N68").Value = -7572.875
$ws.Range("K68").Value = 1824.6666
$ws.Range("H68").Value = 4253.357
$ws.Range("L68").Value = 6074.875
$ws.Range("I71").Value = 1824.6666
$ws.Range("J71").Value = 6074.875
$ws.Range("M71").Value = -5379.333000000001
$ws.Range("L71").Value = 30374.375
$ws.Range("H71").Value = 4253.357
$ws.Range("K71").Value = 9123.333
$ws.Range("N71").Value = -37862.375
$ws.Range("M82").Value = -1516.0
$ws.Range("K82").Value = 1877.0
$ws.Range("I82").Value = 1877.0
$ws.Range("H82").Value = 4791.8335
$ws.Range("K85").Value = 1877.0
$ws.Range("H85").Value = 4791.8335
$ws.Range("I85").Value = 1877.0
$ws.Range("M85").Value = -629.0
$ws.Range("L109").Value = 68997.5
$ws.Range("H109").Value = 68997.5
$ws.Range("N109").Value = -71771.5
$ws.Range("J109").Value = 68997.5
$ws.Range("H132").Value = 4996.161
$ws.Range("K132").Value = 13324.08
$ws.Range("M132").Value = -10794.08
$ws.Range("I132").Value = 4441.36
$ws.Range("L136").Value = 0.0
$ws.Range("H136").Value = 11829.8
$ws.Range("I136").Value = 11829.8
$ws.Range("K136").Value = 35489.39999999999
$ws.Range("J136").Value = 0.0
$ws.Range("N136").Value = $null
$ws.Range("M136").Value = -32939.39999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L81").Value = 9886.667
$ws.Range("K81").Value = 2352.2666
$ws.Range("H81").Value = 2588.8333
$ws.Range("I81").Value = 1176.1333
$ws.Range("J81").Value = 4943.3335
$ws.Range("N81").Value = -12008.667
$ws.Range("M81").Value = -1291.2666
$ws.Range("N84").Value = -60041.335
$ws.Range("J84").Value = 4943.3335
$ws.Range("H84").Value = 2588.8333
$ws.Range("M84").Value = -6457.332999999999
$ws.Range("I84").Value = 1176.1333
$ws.Range("K84").Value = 11761.333
$ws.Range("L84").Value = 49433.335
$ws.Range("H132").Value = 1999.95
$ws.Range("K132").Value = 5468.5386
$ws.Range("M132").Value = -2938.5386
$ws.Range("I132").Value = 1822.8462
$ws.Range("H136").Value = 5271.0967
$ws.Range("I136").Value = 5300.2144
$ws.Range("K136").Value = 15900.6432
$ws.Range("M136").Value = -13350.6432
